$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 589
$ws.Range("F3").Value = 277
$ws.Range("F4").Value = 624
$ws.Range("F6").Value = 782
$ws.Range("F7").Value = 408
$ws.Range("F8").Value = 854
$ws.Range("F9").Value = 465
$ws.Range("F10").Value = 7036
$ws.Range("F11").Value = 1929
$ws.Range("F12").Value = 5120
$ws.Range("F13").Value = 503
$ws.Range("F15").Value = 6648
$ws.Range("F16").Value = 8265
$ws.Range("F18").Value = 1118
$ws.Range("F19").Value = 814
$ws.Range("F20").Value = 4241
$ws.Range("F21").Value = 632
$ws.Range("F22").Value = 121
$ws.Range("F23").Value = 77
$ws.Range("F25").Value = 153
$ws.Range("F26").Value = 1139
$ws.Range("F27").Value = 56
$ws.Range("F28").Value = 1567
$ws.Range("F29").Value = 645
$ws.Range("F30").Value = 793
$ws.Range("F31").Value = 1791
$ws.Range("F32").Value = 295
$ws.Range("F33").Value = 2122
$ws.Range("F34").Value = 275
$ws.Range("F36").Value = 1366
$ws.Range("F38").Value = 758
$ws.Range("F39").Value = 365
$ws.Range("F40").Value = 2842
$ws.Range("F41").Value = 3898
$ws.Range("F43").Value = 27
$ws.Range("F44").Value = 388
$ws.Range("F45").Value = 485
$ws.Range("F46").Value = 9
$ws.Range("F48").Value = 138
$ws.Range("F49").Value = 4024

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1345
$ws.Range("F17").Value = 21
$ws.Range("F19").Value = 32
$ws.Range("F32").Value = 30

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4877

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4877
$ws.Range("F4").Value = 589
$ws.Range("F5").Value = 1345
$ws.Range("F8").Value = 277
$ws.Range("F9").Value = 624
$ws.Range("F11").Value = 782
$ws.Range("F12").Value = 408
$ws.Range("F13").Value = 854
$ws.Range("F14").Value = 465
$ws.Range("F17").Value = 5120
$ws.Range("F18").Value = 6648
$ws.Range("F19").Value = 6648
$ws.Range("F20").Value = 21
$ws.Range("F22").Value = 1118
$ws.Range("F23").Value = 814
$ws.Range("F24").Value = 4241
$ws.Range("F25").Value = 632
$ws.Range("F26").Value = 121
$ws.Range("F28").Value = 153
$ws.Range("F29").Value = 1139
$ws.Range("F30").Value = 1567
$ws.Range("F31").Value = 645
$ws.Range("F32").Value = 793
$ws.Range("F33").Value = 1791
$ws.Range("F34").Value = 295
$ws.Range("F35").Value = 2122
$ws.Range("F40").Value = 758
$ws.Range("F42").Value = 365
$ws.Range("F44").Value = 3898
$ws.Range("F45").Value = 27
$ws.Range("F46").Value = 388
$ws.Range("F48").Value = 138
$ws.Range("F50").Value = 4024
